$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$volCell = $ws.Range("A8")
$volCell.Characters(21,2).Text = "12"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27,9).Text = "3/17/2025"
$weekCell.Characters(47,9).Text = "3/23/2025"

# --- Helper functions for value/type/style changes in the crime stats table ---
function Set-Num($addr, $num) {
    $ws.Range($addr).Value = $num
}

function Set-Text($addr, $text, $styleDonor) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($styleDonor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

function Set-NumRestyle($addr, $num, $styleDonor) {
    $ws.Range($addr).Value = $num
    $ws.Range($styleDonor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# --- Row 15-28 data updates (new crime data collected) ---
Set-Text "C15" "0" "D15"
Set-Num "M15" 50
Set-Num "N15" 20
Set-Num "C16" 1
Set-Num "D16" 2
Set-Num "E16" -50
Set-Num "F16" 8
Set-Num "G16" 13
Set-Num "H16" -38.461538461538
Set-Num "I16" 27
Set-Num "J16" 30
Set-Num "K16" -10
Set-Num "L16" 12.5
Set-Num "M16" -20.588235294117
Set-Num "N16" -85.164835164835
Set-Num "C17" 1
Set-Num "D17" 6
Set-Num "E17" -83.333333333333
Set-Num "F17" 8
Set-Num "G17" 12
Set-Num "H17" -33.333333333333
Set-Num "I17" 27
Set-Num "J17" 36
Set-Num "K17" -25
Set-Num "L17" -38.636363636363
Set-Num "M17" 17.391304347826
Set-Num "N17" -54.237288135593
Set-NumRestyle "C18" 1 "F15"
Set-Num "G18" 1
Set-Num "H18" 200
Set-Num "I18" 14
Set-Num "K18" 16.666666666666
Set-Num "L18" -36.363636363636
Set-Num "M18" 0
Set-Num "N18" -90.140845070422
Set-Num "C19" 6
Set-Num "D19" 6
Set-Num "E19" 0
Set-Num "F19" 30
Set-Num "G19" 26
Set-Num "H19" 15.384615384615
Set-Num "I19" 62
Set-Num "J19" 82
Set-Num "K19" -24.390243902439
Set-Num "L19" -28.735632183908
Set-Num "M19" 16.981132075471
Set-Num "N19" -59.210526315789
Set-Num "C20" 3
Set-Num "F20" 7
Set-Num "G20" 2
Set-Num "H20" 250
Set-Num "I20" 8
Set-Num "K20" 33.333333333333
Set-Num "L20" -61.904761904761
Set-Num "M20" 300
Set-Num "N20" -92.307692307692
Set-Num "C21" 12
Set-Num "D21" 14
Set-Num "E21" -14.285714285714
Set-Num "F21" 58
Set-Num "G21" 54
Set-Num "H21" 7.407407407407
Set-Num "I21" 144
Set-Num "J21" 168
Set-Num "L21" -28
Set-Num "M21" 9.090909090909
Set-Num "N21" -77.777777777777
Set-Text "D22" "0" "D15"
Set-Text "E22" "***.*" "D15"
Set-Num "I23" 26
Set-Num "J23" 35
Set-Num "K23" -25.714285714285
Set-Num "L23" -23.529411764705
Set-Num "M23" 23.809523809523
Set-Num "C24" 8
Set-Num "D24" 7
Set-Num "E24" 14.285714285714
Set-Num "F24" 39
Set-Num "G24" 34
Set-Num "H24" 14.705882352941
Set-Num "I24" 110
Set-Num "J24" 85
Set-Num "K24" 29.411764705882
Set-Num "L24" 7.843137254901
Set-Num "M24" 22.222222222222
Set-NumRestyle "D25" 1 "F15"
Set-NumRestyle "E25" 100 "K15"
Set-Num "F25" 7
Set-Num "G25" 3
Set-Num "H25" 133.333333333333
Set-Num "I25" 17
Set-Num "J25" 15
Set-Num "K25" 13.333333333333
Set-Num "L25" -52.777777777777
Set-Num "C26" 4
Set-Num "D26" 5
Set-Num "E26" -20
Set-Num "F26" 24
Set-Num "G26" 21
Set-Num "H26" 14.285714285714
Set-Num "I26" 59
Set-Num "J26" 64
Set-Num "K26" -7.8125
Set-Num "L26" 7.272727272727
Set-Num "M26" 9.259259259259
Set-Text "C27" "0" "D15"
Set-Text "C28" "0" "D15"
Set-NumRestyle "D28" 1 "F15"
Set-NumRestyle "E28" -100 "K15"
Set-Num "G28" 2
Set-Num "H28" 0
Set-Num "J28" 4
Set-Num "K28" 50
